# "Version 1." -> "Version 2." but authored the way a real Word edit
# session would leave it behind: the word "Version" ends up split into
# two runs ("Versi" + "on") around the spell-check markers, " 1." loses
# its trailing period, and a lone "." run reappears after the _GoBack
# bookmark.

$d = $word.ActiveDocument

# --- Locate "Version" and split it into "Versi" | "on" -------------------
# Word only splits a run when something actually changes at the split
# point (a format toggle, a bookmark, etc). Dropping a throwaway bookmark
# in the middle of the run and immediately deleting it is the cleanest
# way to get a bare run boundary with no leftover direct formatting.
$findVersion = $d.Content
[void]$findVersion.Find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($findVersion.Start + 5, $findVersion.Start + 5)
$tempBookmarkName = "__split_marker__"
[void]$d.Bookmarks.Add($tempBookmarkName, $splitPoint)
[void]$d.Bookmarks($tempBookmarkName).Delete()

# --- " 1." -> " 2" (drop the trailing period from this run) --------------
$findNumber = $d.Content
[void]$findNumber.Find.Execute(" 1.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findNumber.Text = " 2"

# --- Re-append the final "." as its own run after the _GoBack bookmark ---
$goBack = $d.Bookmarks("_GoBack")
[void]$goBack.Range.InsertAfter(".")
